# Update workbook for data through 2022-01-03 (adds a day of data, 12-25 -> 12-26)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet (tab) name
$ws.Name = "Through 2021-12-26"

# Update the "December" row label to reflect new through-date
$ws.Range("A13").Value = "December (through 12-26)"

# Update December row (row 13) values
$ws.Range("B13").Value = 39
$ws.Range("C13").Value = 85
$ws.Range("D13").Value = 101
$ws.Range("E13").Value = 59
$ws.Range("F13").Value = 55
$ws.Range("G13").Value = 121
$ws.Range("H13").Value = 161

# Update Total row (row 14) values
$ws.Range("B14").Value = 330
$ws.Range("C14").Value = 648
$ws.Range("D14").Value = 922
$ws.Range("E14").Value = 741
$ws.Range("F14").Value = 589
$ws.Range("G14").Value = 1385
$ws.Range("H14").Value = 1804
